$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2367.4375
$ws.Range("I40").Value = 3545.8
$ws.Range("K40").Value = 3545.8
$ws.Range("M40").Value = -3370.8
$ws.Range("H125").Value = 55556556
$ws.Range("I125").Value = 142857470
$ws.Range("J125").Value = 1431.8182
$ws.Range("K125").Value = 1285717230
$ws.Range("L125").Value = 12886.3638
$ws.Range("M125").Value = -1285714770
$ws.Range("N125").Value = -17806.3638
$ws.Range("H132").Value = 2978132.8
$ws.Range("I132").Value = 3761315.2
$ws.Range("J132").Value = 2038.4
$ws.Range("K132").Value = 11283945.6
$ws.Range("L132").Value = 6115.200000000001
$ws.Range("M132").Value = -11281415.6
$ws.Range("N132").Value = -11175.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2866.44
$ws.Range("I32").Value = 2838.8542
$ws.Range("J32").Value = 3528.5
$ws.Range("K32").Value = 2838.8542
$ws.Range("L32").Value = 3528.5
$ws.Range("M32").Value = -2551.8542
$ws.Range("N32").Value = -4102.5
$ws.Range("H61").Value = 1829.7407
$ws.Range("I61").Value = 1469.421
$ws.Range("J61").Value = 2685.5
$ws.Range("K61").Value = 1469.421
$ws.Range("L61").Value = 2685.5
$ws.Range("M61").Value = -1257.421
$ws.Range("N61").Value = -3109.5
$ws.Range("H74").Value = 975.6585
$ws.Range("I74").Value = 975.05
$ws.Range("K74").Value = 975.05
$ws.Range("M74").Value = -101.05
$ws.Range("H77").Value = 975.6585
$ws.Range("I77").Value = 975.05
$ws.Range("K77").Value = 4875.25
$ws.Range("M77").Value = -507.25
$ws.Range("H122").Value = 1590.1714
$ws.Range("I122").Value = 1518.8572
$ws.Range("J122").Value = 1875.4286
$ws.Range("K122").Value = 4556.571599999999
$ws.Range("L122").Value = 5626.2858
$ws.Range("M122").Value = -2106.571599999999
$ws.Range("N122").Value = -10526.2858
$ws.Range("H132").Value = 3141.986
$ws.Range("I132").Value = 3100.8245
$ws.Range("J132").Value = 3298.4
$ws.Range("K132").Value = 9302.4735
$ws.Range("L132").Value = 9895.200000000001
$ws.Range("M132").Value = -6772.4735
$ws.Range("N132").Value = -14955.2
$ws.Range("H136").Value = 1829.7407
$ws.Range("I136").Value = 1469.421
$ws.Range("J136").Value = 2685.5
$ws.Range("K136").Value = 4408.263
$ws.Range("L136").Value = 8056.5
$ws.Range("M136").Value = -1858.263
$ws.Range("N136").Value = -13156.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4239.9165
$ws.Range("I20").Value = 4737.16
$ws.Range("J20").Value = 3109.818
$ws.Range("K20").Value = 4737.16
$ws.Range("L20").Value = 3109.818
$ws.Range("M20").Value = -4490.16
$ws.Range("N20").Value = -3603.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2690969.5
$ws.Range("I31").Value = 2143.6584
$ws.Range("J31").Value = 7940581.5
$ws.Range("K31").Value = 2143.6584
$ws.Range("L31").Value = 7940581.5
$ws.Range("M31").Value = -1848.6584
$ws.Range("N31").Value = -7941171.5
$ws.Range("H34").Value = 2690969.5
$ws.Range("I34").Value = 2143.6584
$ws.Range("J34").Value = 7940581.5
$ws.Range("K34").Value = 2143.6584
$ws.Range("L34").Value = 7940581.5
$ws.Range("M34").Value = -1941.6584
$ws.Range("N34").Value = -7940985.5
$ws.Range("H58").Value = 2143.5
$ws.Range("I58").Value = 2274.7273
$ws.Range("J58").Value = 700
$ws.Range("K58").Value = 2274.7273
$ws.Range("L58").Value = 700
$ws.Range("M58").Value = -2071.7273
$ws.Range("N58").Value = -1106
$ws.Range("H94").Value = 111112240
$ws.Range("J94").Value = 1265.75
$ws.Range("L94").Value = 1265.75
$ws.Range("N94").Value = -2167.75
$ws.Range("H105").Value = 996.55554
$ws.Range("I105").Value = 928
$ws.Range("J105").Value = 1133.6666
$ws.Range("K105").Value = 928
$ws.Range("L105").Value = 1133.6666
$ws.Range("M105").Value = 819
$ws.Range("N105").Value = -4627.6666
$ws.Range("H132").Value = 5003474
$ws.Range("I132").Value = 3805.5
$ws.Range("J132").Value = 7356259.5
$ws.Range("K132").Value = 11416.5
$ws.Range("L132").Value = 22068778.5
$ws.Range("M132").Value = -8886.5
$ws.Range("N132").Value = -22073838.5
$ws.Range("H136").Value = 2143.5
$ws.Range("I136").Value = 2274.7273
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 6824.1819
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = -4274.1819
$ws.Range("N136").Value = -7200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 15600
$ws.Range("I119").Value = 1000
$ws.Range("J119").Value = 17222.223
$ws.Range("K119").Value = 3000
$ws.Range("L119").Value = 51666.66900000001
$ws.Range("M119").Value = 1838
$ws.Range("N119").Value = -61342.66900000001
$ws.Range("H120").Value = 23972.37
$ws.Range("I120").Value = 13670
$ws.Range("J120").Value = 25904.062
$ws.Range("K120").Value = 41010
$ws.Range("L120").Value = 77712.186
$ws.Range("M120").Value = -36172
$ws.Range("N120").Value = -87388.186
$ws.Range("H131").Value = 744.91
$ws.Range("I131").Value = 292.5
$ws.Range("J131").Value = 806.6023
$ws.Range("K131").Value = 877.5
$ws.Range("L131").Value = 2419.8069
$ws.Range("M131").Value = 4162.5
$ws.Range("N131").Value = -12499.8069

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 38463930
$ws.Range("I122").Value = 76925464
$ws.Range("J122").Value = 2390.5386
$ws.Range("K122").Value = 230776392
$ws.Range("L122").Value = 7171.6158
$ws.Range("M122").Value = -230773942
$ws.Range("N122").Value = -12071.6158
$ws.Range("H123").Value = 26490
$ws.Range("J123").Value = 26490
$ws.Range("L123").Value = 26490
$ws.Range("N123").Value = -31390
$ws.Range("H126").Value = 2619.6428
$ws.Range("I126").Value = 3119.4443
$ws.Range("K126").Value = 9358.332900000001
$ws.Range("M126").Value = -6888.332900000001
$ws.Range("H131").Value = 21000
$ws.Range("J131").Value = 21000
$ws.Range("L131").Value = 21000
$ws.Range("N131").Value = -31080
$ws.Range("H132").Value = 2214.7932
$ws.Range("I132").Value = 1726.238
$ws.Range("J132").Value = 3497.25
$ws.Range("K132").Value = 5178.714
$ws.Range("L132").Value = 10491.75
$ws.Range("M132").Value = -2648.714
$ws.Range("N132").Value = -15551.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1557.8
$ws.Range("I16").Value = 1447.25
$ws.Range("K16").Value = 1447.25
$ws.Range("M16").Value = -1277.25
$ws.Range("H22").Value = 487.69232
$ws.Range("I22").Value = 644.5
$ws.Range("J22").Value = 353.2857
$ws.Range("K22").Value = 644.5
$ws.Range("L22").Value = 353.2857
$ws.Range("M22").Value = -349.5
$ws.Range("N22").Value = -943.2857
$ws.Range("H27").Value = 487.69232
$ws.Range("I27").Value = 644.5
$ws.Range("J27").Value = 353.2857
$ws.Range("K27").Value = 644.5
$ws.Range("L27").Value = 353.2857
$ws.Range("M27").Value = -537.5
$ws.Range("N27").Value = -567.2857
$ws.Range("H132").Value = 8030.2812
$ws.Range("I132").Value = 12204.667
$ws.Range("J132").Value = 2663.2144
$ws.Range("K132").Value = 36614.001
$ws.Range("L132").Value = 7989.6432
$ws.Range("M132").Value = -34084.001
$ws.Range("N132").Value = -13049.6432
$ws.Range("H133").Value = 26000
$ws.Range("J133").Value = 26000
$ws.Range("L133").Value = 26000
$ws.Range("N133").Value = -31060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19928000
$ws.Range("I107").Value = 8333812.5
$ws.Range("K107").Value = 25001437.5
$ws.Range("M107").Value = -24999517.5
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H126").Value = 52639250
$ws.Range("I126").Value = 52639250
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 157917750
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -157915280
$ws.Range("N110").ClearContents()
$ws.Range("N126").ClearContents()
